$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-12-21 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-12-22 Sunday", 2)
$d.Content.Find.Execute("421÷8=52, 5", $true, $false, $false, $false, $false, $true, 1, $false, "496÷4=124, 0", 2)
$d.Content.Find.Execute("847÷6=141, 1", $true, $false, $false, $false, $false, $true, 1, $false, "677÷4=169, 1", 2)
$d.Content.Find.Execute("535÷8=66, 7", $true, $false, $false, $false, $false, $true, 1, $false, "740÷6=123, 2", 2)
$d.Content.Find.Execute("776÷3=258, 2", $true, $false, $false, $false, $false, $true, 1, $false, "288÷7=41, 1", 2)
$d.Content.Find.Execute("988÷8=123, 4", $true, $false, $false, $false, $false, $true, 1, $false, "728÷6=121, 2", 2)
$d.Content.Find.Execute("977÷6=162, 5", $true, $false, $false, $false, $false, $true, 1, $false, "365÷7=52, 1", 2)
$d.Content.Find.Execute("160÷8=20, 0", $true, $false, $false, $false, $false, $true, 1, $false, "534÷5=106, 4", 2)
$d.Content.Find.Execute("638÷4=159, 2", $true, $false, $false, $false, $false, $true, 1, $false, "928÷3=309, 1", 2)
$d.Content.Find.Execute("104÷8=13, 0", $true, $false, $false, $false, $false, $true, 1, $false, "733÷3=244, 1", 2)
$d.Content.Find.Execute("279÷8=34, 7", $true, $false, $false, $false, $false, $true, 1, $false, "715÷4=178, 3", 2)
$d.Content.Find.Execute("684÷4=171, 0", $true, $false, $false, $false, $false, $true, 1, $false, "638÷8=79, 6", 2)
$d.Content.Find.Execute("559÷6=93, 1", $true, $false, $false, $false, $false, $true, 1, $false, "797÷4=199, 1", 2)
$d.Content.Find.Execute("834÷3=278, 0", $true, $false, $false, $false, $false, $true, 1, $false, "202÷5=40, 2", 2)
$d.Content.Find.Execute("445÷2=222, 1", $true, $false, $false, $false, $false, $true, 1, $false, "132÷2=66, 0", 2)
$d.Content.Find.Execute("991÷9=110, 1", $true, $false, $false, $false, $false, $true, 1, $false, "474÷7=67, 5", 2)
$d.Content.Find.Execute("609÷5=121, 4", $true, $false, $false, $false, $false, $true, 1, $false, "809÷4=202, 1", 2)
$d.Content.Find.Execute("610÷3=203, 1", $true, $false, $false, $false, $false, $true, 1, $false, "127÷4=31, 3", 2)
$d.Content.Find.Execute("252÷5=50, 2", $true, $false, $false, $false, $false, $true, 1, $false, "701÷7=100, 1", 2)
$d.Content.Find.Execute("598÷7=85, 3", $true, $false, $false, $false, $false, $true, 1, $false, "941÷4=235, 1", 2)
$d.Content.Find.Execute("180÷4=45, 0", $true, $false, $false, $false, $false, $true, 1, $false, "880÷2=440, 0", 2)
$d.Content.Find.Execute("517÷9=57, 4", $true, $false, $false, $false, $false, $true, 1, $false, "554÷3=184, 2", 2)
$d.Content.Find.Execute("957÷3=319, 0", $true, $false, $false, $false, $false, $true, 1, $false, "384÷6=64, 0", 2)
$d.Content.Find.Execute("347÷4=86, 3", $true, $false, $false, $false, $false, $true, 1, $false, "555÷8=69, 3", 2)
$d.Content.Find.Execute("483÷2=241, 1", $true, $false, $false, $false, $false, $true, 1, $false, "145÷3=48, 1", 2)
$d.Content.Find.Execute("742÷6=123, 4", $true, $false, $false, $false, $false, $true, 1, $false, "554÷2=277, 0", 2)
